$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared-string runs) ---
# A8: "Volume 31   Number  50" -> "...51" (Volume/Number header)
$ws.Range("A8").Value = "Volume 31   Number  51"
# C9: week-covering date range shifts forward by one week
$ws.Range("C9").Value = "Report Covering the Week  12/16/2024  Through  12/22/2024"

# --- CompStat table value updates (rows 14-33) ---
# Row 14
$ws.Range("C14").Value = 1
$ws.Range("F14").Value = 1
$ws.Range("G14").Value = 1
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 17
$ws.Range("K14").Value = -29.166666666666
$ws.Range("L14").Value = -37.037037037037
$ws.Range("M14").Value = -43.333333333333
$ws.Range("N14").Value = -86.290322580645

# Row 15
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 4
$ws.Range("E15").Value = 75
$ws.Range("F15").Value = 20
$ws.Range("G15").Value = 14
$ws.Range("H15").Value = 42.857142857142
$ws.Range("I15").Value = 235
$ws.Range("J15").Value = 194
$ws.Range("K15").Value = 21.134020618556
$ws.Range("L15").Value = 15.763546798029
$ws.Range("M15").Value = 69.064748201438
$ws.Range("N15").Value = 16.336633663366

# Row 16
$ws.Range("C16").Value = 27
$ws.Range("D16").Value = 50
$ws.Range("E16").Value = -46
$ws.Range("F16").Value = 124
$ws.Range("G16").Value = 160
$ws.Range("H16").Value = -22.5
$ws.Range("I16").Value = 2088
$ws.Range("J16").Value = 2010
$ws.Range("K16").Value = 3.880597014925
$ws.Range("L16").Value = 13.725490196078
$ws.Range("M16").Value = 11.36
$ws.Range("N16").Value = -75.827737902292

# Row 17
$ws.Range("C17").Value = 46
$ws.Range("D17").Value = 57
$ws.Range("E17").Value = -19.298245614035
$ws.Range("F17").Value = 218
$ws.Range("G17").Value = 210
$ws.Range("H17").Value = 3.809523809523
$ws.Range("I17").Value = 3212
$ws.Range("J17").Value = 2844
$ws.Range("K17").Value = 12.939521800281
$ws.Range("L17").Value = 32.344458178821
$ws.Range("M17").Value = 112.153236459709
$ws.Range("N17").Value = 12.701754385964

# Row 18
$ws.Range("C18").Value = 58
$ws.Range("D18").Value = 50
$ws.Range("E18").Value = 16
$ws.Range("F18").Value = 191
$ws.Range("G18").Value = 180
$ws.Range("H18").Value = 6.111111111111
$ws.Range("I18").Value = 1965
$ws.Range("J18").Value = 1950
$ws.Range("K18").Value = 0.76923076923
$ws.Range("L18").Value = -3.39233038348
$ws.Range("M18").Value = -24.510180560891
$ws.Range("N18").Value = -86.407028223575

# Row 19
$ws.Range("C19").Value = 116
$ws.Range("D19").Value = 123
$ws.Range("E19").Value = -5.691056910569
$ws.Range("F19").Value = 427
$ws.Range("G19").Value = 521
$ws.Range("H19").Value = -18.042226487524
$ws.Range("I19").Value = 6418
$ws.Range("J19").Value = 6620
$ws.Range("K19").Value = -3.051359516616
$ws.Range("L19").Value = -6.443148688046
$ws.Range("M19").Value = 61.377923057581
$ws.Range("N19").Value = -21.874619598295

# Row 20
$ws.Range("C20").Value = 23
$ws.Range("D20").Value = 53
$ws.Range("E20").Value = -56.603773584905
$ws.Range("F20").Value = 116
$ws.Range("G20").Value = 189
$ws.Range("H20").Value = -38.624338624338
$ws.Range("I20").Value = 2488
$ws.Range("J20").Value = 2622
$ws.Range("K20").Value = -5.11060259344
$ws.Range("L20").Value = 26.809378185525
$ws.Range("M20").Value = 46.009389671361
$ws.Range("N20").Value = -89.225239270711

# Row 21
$ws.Range("C21").Value = 278
$ws.Range("D21").Value = 337
$ws.Range("E21").Value = -17.507418397626
$ws.Range("F21").Value = 1097
$ws.Range("G21").Value = 1275
$ws.Range("H21").Value = -13.960784313725
$ws.Range("I21").Value = 16423
$ws.Range("J21").Value = 16264
$ws.Range("K21").Value = 0.977619281849
$ws.Range("L21").Value = 6.997198514561
$ws.Range("M21").Value = 38.684343860834
$ws.Range("N21").Value = -71.475962206474

# Row 22
$ws.Range("C22").Value = 8
$ws.Range("D22").Value = 2
$ws.Range("E22").Value = 300
$ws.Range("G22").Value = 24
$ws.Range("H22").Value = 16.666666666666
$ws.Range("I22").Value = 321
$ws.Range("J22").Value = 353
$ws.Range("K22").Value = -9.065155807365
$ws.Range("L22").Value = 7.357859531772
$ws.Range("M22").Value = 69.841269841269

# Row 23
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 5
$ws.Range("E23").Value = -60
$ws.Range("F23").Value = 15
$ws.Range("G23").Value = 18
$ws.Range("H23").Value = -16.666666666666
$ws.Range("I23").Value = 241
$ws.Range("J23").Value = 261
$ws.Range("K23").Value = -7.662835249042
$ws.Range("L23").Value = 2.118644067796
$ws.Range("M23").Value = 43.45238095238

# Row 24
$ws.Range("C24").Value = 287
$ws.Range("D24").Value = 341
$ws.Range("E24").Value = -15.835777126099
$ws.Range("F24").Value = 1175
$ws.Range("G24").Value = 1366
$ws.Range("H24").Value = -13.982430453879
$ws.Range("I24").Value = 15520
$ws.Range("J24").Value = 15276
$ws.Range("K24").Value = 1.597276774024
$ws.Range("L24").Value = 2.672664726118
$ws.Range("M24").Value = 66.863778088377

# Row 25
$ws.Range("C25").Value = 162
$ws.Range("D25").Value = 201
$ws.Range("E25").Value = -19.402985074626
$ws.Range("F25").Value = 672
$ws.Range("G25").Value = 812
$ws.Range("H25").Value = -17.241379310344
$ws.Range("I25").Value = 9471
$ws.Range("J25").Value = 8626
$ws.Range("K25").Value = 9.795965685137
$ws.Range("L25").Value = 21.190019193858

# Row 26
$ws.Range("C26").Value = 91
$ws.Range("D26").Value = 110
$ws.Range("E26").Value = -17.272727272727
$ws.Range("F26").Value = 425
$ws.Range("G26").Value = 398
$ws.Range("H26").Value = 6.783919597989
$ws.Range("I26").Value = 6113
$ws.Range("J26").Value = 5296
$ws.Range("K26").Value = 15.42673716012
$ws.Range("L26").Value = 26.484585143803
$ws.Range("M26").Value = 31.434100193506

# Row 27
$ws.Range("C27").Value = 7
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 75
$ws.Range("F27").Value = 21
$ws.Range("G27").Value = 21
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 326
$ws.Range("J27").Value = 304
$ws.Range("K27").Value = 7.236842105263
$ws.Range("L27").Value = 13.194444444444

# Row 28
$ws.Range("C28").Value = 13
$ws.Range("D28").Value = 9
$ws.Range("E28").Value = 44.444444444444
$ws.Range("F28").Value = 38
$ws.Range("H28").Value = -17.391304347826
$ws.Range("I28").Value = 633
$ws.Range("J28").Value = 680
$ws.Range("K28").Value = -6.911764705882
$ws.Range("L28").Value = 1.605136436597

# Row 29
$ws.Range("C29").Value = 5
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Value = "'0"
$ws.Range("E29").Value = "***.*"
$ws.Range("F29").Value = 8
$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 100
$ws.Range("I29").Value = 33
$ws.Range("K29").Value = -50
$ws.Range("L29").Value = -57.142857142857
$ws.Range("M29").Value = -28.260869565217
$ws.Range("N29").Value = -86.852589641434

# Row 30
$ws.Range("C30").Value = 3
$ws.Range("D30").NumberFormat = "General"
$ws.Range("D30").Value = "'0"
$ws.Range("E30").Value = "***.*"
$ws.Range("F30").Value = 5
$ws.Range("G30").Value = 3
$ws.Range("H30").Value = 66.666666666666
$ws.Range("I30").Value = 27
$ws.Range("K30").Value = -55
$ws.Range("L30").Value = -57.8125
$ws.Range("M30").Value = -28.947368421052
$ws.Range("N30").Value = -88.053097345132

# Row 31
$ws.Range("D31").NumberFormat = "General"
$ws.Range("D31").Value = "'0"
$ws.Range("E31").Value = "***.*"
$ws.Range("F31").Value = 4
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = -20
$ws.Range("I31").Value = 73
$ws.Range("K31").Value = -5.194805194805
$ws.Range("L31").Value = 23.728813559322

# Row 33
$ws.Range("C33").Value = 1
$ws.Range("D33").Value = 1
$ws.Range("E33").Value = 0
$ws.Range("I33").Value = 45
$ws.Range("J33").Value = 42
$ws.Range("K33").Value = 7.142857142857
$ws.Range("L33").Value = 18.421052631578

# --- Insert new blank row before row 56 (shifts old 56/57 -> 57/58) ---
$ws.Rows.Item(56).Insert()
